# Automatische test-sync: 2025-07-29 21:53:50
# Appends a new mail-log row (#12) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover it, and swaps the "Bestelling /
# Levering" / "Retour / Terugbetaling" rows on the "Dashboard" sheet
# (with the updated counts) to match the refreshed aggregate.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: add row 14
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Ik heb nog geen geld terug."
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Testmail #12: Ik heb nog geen geld terug."
$logs.Range("D14").Value = "Retour / Terugbetaling"
$logs.Range("E14").Value = "Beste klant,`nBedankt voor uw e-mail. Om uw vraag beter te kunnen begrijpen en u verder te kunnen helpen, zou ik graag wat meer informatie ontvangen. Kunt u mij uw ordernummer of referentienummer doorgeven, zodat ik kan controleren wat de status is van uw terugbetaling?`nIk kijk uit naar uw reactie.`nMet vriendelijke groet,`n[Naam] E-mailassistent `n[Bedrijfsnaam]"
$logs.Range("F14").Value = "2025-07-29 21:53:18"
$logs.Range("G14").Value = "Ja"
$logs.Range("H14").Value = "Nee"
$logs.Range("I14").Value = "Ja"
$logs.Range("J14").Value = "Nee"

# Writing the multi-line reply text auto-expands the row height; the
# other multi-line "Antwoord" rows (6, 7, 13, ...) keep the default
# height, so put row 14 back in line with them.
$logs.Rows.Item(14).AutoFit()

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 13 to
# row 14 so the new row is covered, same as the original authoring tool
# would do when the sqref grows with the data.
$logs.Range("D2:D13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D14"))
$logs.Range("G2:G13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G14"))
$logs.Range("H2:H13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H14"))
$logs.Range("I2:I13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I14"))
$logs.Range("J2:J13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J14"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: the new "Retour / Terugbetaling" mail bumps its
#    category count to 2, pushing it above "Bestelling / Levering" (1)
#    in the summary table, so rows 5 & 6 swap places.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Retour / Terugbetaling"
$dashboard.Range("B5").Value = 2
$dashboard.Range("A6").Value = "Bestelling / Levering"
$dashboard.Range("B6").Value = 1
